# chore: update Sheets via scheduled runner
# Refresh cached market-board figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the leve-profit tracker sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8186.75
$ws.Range("I32").Value = 5748.5
$ws.Range("J32").Value = 8999.5
$ws.Range("K32").Value = 5748.5
$ws.Range("L32").Value = 8999.5
$ws.Range("M32").Value = -5422.5
$ws.Range("N32").Value = -9651.5
$ws.Range("H137").Value = 2500
$ws.Range("I137").Value = 2500
$ws.Range("K137").Value = 7500
$ws.Range("M137").Value = -4950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 20288
$ws.Range("J76").Value = 20288
$ws.Range("L76").Value = 20288
$ws.Range("N76").Value = -20964
$ws.Range("H79").Value = 20288
$ws.Range("J79").Value = 20288
$ws.Range("L79").Value = 20288
$ws.Range("N79").Value = -22628
$ws.Range("H132").Value = 812
$ws.Range("I132").Value = 812
$ws.Range("K132").Value = 2436
$ws.Range("M132").Value = 94
$ws.Range("H135").Value = 61497
$ws.Range("J135").Value = 61497
$ws.Range("L135").Value = 61497
$ws.Range("N135").Value = -71637

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3995
$ws.Range("J31").Value = 4863.4287
$ws.Range("L31").Value = 4863.4287
$ws.Range("N31").Value = -5453.4287
$ws.Range("H33").Value = 1248.2
$ws.Range("I33").Value = 1248.2
$ws.Range("K33").Value = 1248.2
$ws.Range("M33").Value = -869.2
$ws.Range("H34").Value = 3995
$ws.Range("J34").Value = 4863.4287
$ws.Range("L34").Value = 4863.4287
$ws.Range("N34").Value = -5267.4287
$ws.Range("H36").Value = 3548
$ws.Range("I36").Value = 3548
$ws.Range("K36").Value = 3548
$ws.Range("M36").Value = -3160
$ws.Range("H40").Value = 3548
$ws.Range("I40").Value = 3548
$ws.Range("K40").Value = 3548
$ws.Range("M40").Value = -3388
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 999
$ws.Range("I59").Value = 999
$ws.Range("J59").Value = 999
$ws.Range("K59").Value = 2997
$ws.Range("L59").Value = 2997
$ws.Range("M59").Value = -2457
$ws.Range("N59").Value = -4077

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30820
$ws.Range("H106").Value = 8000
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 8000
$ws.Range("N106").Value = -10524
$ws.Range("H137").Value = 51500
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 5836.1665
$ws.Range("I33").Value = 7754.5
$ws.Range("J33").Value = 1999.5
$ws.Range("K33").Value = 7754.5
$ws.Range("L33").Value = 1999.5
$ws.Range("M33").Value = -7504.5
$ws.Range("N33").Value = -2499.5
$ws.Range("H34").Value = 500
$ws.Range("I34").Value = 500
$ws.Range("K34").Value = 500
$ws.Range("M34").Value = -297
$ws.Range("H36").Value = 5836.1665
$ws.Range("I36").Value = 7754.5
$ws.Range("J36").Value = 1999.5
$ws.Range("K36").Value = 7754.5
$ws.Range("L36").Value = 1999.5
$ws.Range("M36").Value = -7504.5
$ws.Range("N36").Value = -2499.5
$ws.Range("H119").Value = 35000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 35000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H122").Value = 2777.4
$ws.Range("I122").Value = 3222.5
$ws.Range("J122").Value = 997
$ws.Range("K122").Value = 9667.5
$ws.Range("L122").Value = 2991
$ws.Range("M122").Value = -7217.5
$ws.Range("N122").Value = -7891
$ws.Range("H123").Value = 45000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 45000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800
$ws.Range("H124").Value = 47500
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 47500
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 47500
$ws.Range("N124").Value = -57320
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 3562.2222
$ws.Range("I126").Value = 3007.5
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 9022.5
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -6552.5
$ws.Range("N126").Value = -28940
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 99995
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 99995
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 99995
$ws.Range("N129").Value = -109995
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 99995
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 99995
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 99995
$ws.Range("N131").Value = -110075
$ws.Range("H132").Value = 943.1667
$ws.Range("I132").Value = 952.5
$ws.Range("J132").Value = 924.5
$ws.Range("K132").Value = 2857.5
$ws.Range("L132").Value = 2773.5
$ws.Range("M132").Value = -327.5
$ws.Range("N132").Value = -7833.5
$ws.Range("H133").Value = 95000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 95000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 95000
$ws.Range("N133").Value = -105120
$ws.Range("H135").Value = 99464
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 99464
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 99464
$ws.Range("N135").Value = -109604
$ws.Range("H136").Value = 3197.8
$ws.Range("I136").Value = 2247.5
$ws.Range("J136").Value = 6999
$ws.Range("K136").Value = 6742.5
$ws.Range("L136").Value = 20997
$ws.Range("M136").Value = -4192.5
$ws.Range("N136").Value = -26097
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
